# Split the former combined "experiment" workbook's beads/cells sheets
# apart (fixes #45): the "cells" sheet still had a couple of stray
# leftover values from the "beads" sheet (a Gate Fraction in F2 and a
# duplicated Beads File Path in B3) that belong only on the beads sheet.
# Clear them here so "cells" becomes a clean per-file template row set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cells")

# Remove the leftover Gate Fraction value on row 2 (FCFiles/data_001.fcs).
$ws.Range("F2").Clear()

# Remove the leftover/duplicated Beads File Path value on row 3
# (FCFiles/data_002.fcs).
$ws.Range("B3").Clear()

# Leave the cursor where the author ended up after making the edits.
$ws.Range("C22").Select()
